# "Gráficos todos os classificadores"
# Align the last "Branch and Bound" rows (158-169) with the short-form
# naming convention ("Branch and Bound <Balanceado/Desbalanceado> -
# <Normalizado/Nao Normalizado>") already used everywhere else in the
# sheet (e.g. rows 62-73, 98-109, 146-157), instead of the older long-form
# wording ("Branch and Bound - <Normalizado/Sem Normalização> -
# <Desbalanceado/Balanceado>").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C158").Value = "Branch and Bound Desbalanceado - Normalizado"
$ws.Range("C159").Value = "Branch and Bound Desbalanceado - Normalizado"
$ws.Range("C160").Value = "Branch and Bound Desbalanceado - Normalizado"

$ws.Range("C161").Value = "Branch and Bound Balanceado - Normalizado"
$ws.Range("C162").Value = "Branch and Bound Balanceado - Normalizado"
$ws.Range("C163").Value = "Branch and Bound Balanceado - Normalizado"

$ws.Range("C164").Value = "Branch and Bound Desbalanceado - Nao Normalizado"
$ws.Range("C165").Value = "Branch and Bound Desbalanceado - Nao Normalizado"
$ws.Range("C166").Value = "Branch and Bound Desbalanceado - Nao Normalizado"

$ws.Range("C167").Value = "Branch and Bound Balanceado - Nao Normalizado"
$ws.Range("C168").Value = "Branch and Bound Balanceado - Nao Normalizado"
$ws.Range("C169").Value = "Branch and Bound Balanceado - Nao Normalizado"

# Column C now holds the longer short-form labels; widen it (and drop the
# stale autofit width) so the text is not clipped, matching the workbook
# as last saved.
$ws.Columns.Item(3).ColumnWidth = 67.5

# Leave the selection where it was left in the source workbook.
$ws.Range("C166").Select() | Out-Null
